$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 'b'
$ws.Range("J2").Value = 'Acknowledge (Backchannel)'
$ws.Range("I4").Value = '%'
$ws.Range("J4").Value = 'Uninterpretable'
$ws.Range("I6").Value = 'qy'
$ws.Range("J6").Value = 'Yes-No-Question'
$ws.Range("I7").Value = 'ba'
$ws.Range("J7").Value = 'Appreciation'
$ws.Range("I18").Value = 'aa'
$ws.Range("J18").Value = 'Agree/Accept'
$ws.Range("I19").Value = 'ba'
$ws.Range("J19").Value = 'Appreciation'
$ws.Range("I27").Value = 'qy'
$ws.Range("J27").Value = 'Yes-No-Question'
$ws.Range("I32").Value = 'sd'
$ws.Range("J32").Value = 'Statement-non-opinion'
$ws.Range("I33").Value = 'sd'
$ws.Range("J33").Value = 'Statement-non-opinion'
$ws.Range("I44").Value = 'sv'
$ws.Range("J44").Value = 'Statement-opinion'
$ws.Range("I45").Value = 'b'
$ws.Range("J45").Value = 'Acknowledge (Backchannel)'
$ws.Range("I69").Value = 'sd'
$ws.Range("J69").Value = 'Statement-non-opinion'
$ws.Range("I71").Value = 'aa'
$ws.Range("J71").Value = 'Agree/Accept'
$ws.Range("I76").Value = 'sd'
$ws.Range("J76").Value = 'Statement-non-opinion'
$ws.Range("I78").Value = 'aa'
$ws.Range("J78").Value = 'Agree/Accept'
$ws.Range("I84").Value = 'aa'
$ws.Range("J84").Value = 'Agree/Accept'
$ws.Range("I85").Value = 'b'
$ws.Range("J85").Value = 'Acknowledge (Backchannel)'
$ws.Range("I88").Value = '%'
$ws.Range("J88").Value = 'Uninterpretable'
$ws.Range("I89").Value = 'b'
$ws.Range("J89").Value = 'Acknowledge (Backchannel)'
$ws.Range("I90").Value = 'aa'
$ws.Range("J90").Value = 'Agree/Accept'
$ws.Range("I91").Value = 'b'
$ws.Range("J91").Value = 'Acknowledge (Backchannel)'
$ws.Range("I103").Value = 'sd'
$ws.Range("J103").Value = 'Statement-non-opinion'
$ws.Range("I123").Value = 'sv'
$ws.Range("J123").Value = 'Statement-opinion'
$ws.Range("I126").Value = 'sv'
$ws.Range("J126").Value = 'Statement-opinion'
$ws.Range("I132").Value = '%'
$ws.Range("J132").Value = 'Uninterpretable'
$ws.Range("I155").Value = 'aa'
$ws.Range("J155").Value = 'Agree/Accept'
$ws.Range("I161").Value = 'aa'
$ws.Range("J161").Value = 'Agree/Accept'
$ws.Range("I163").Value = 'sd'
$ws.Range("J163").Value = 'Statement-non-opinion'
$ws.Range("I168").Value = 'aa'
$ws.Range("J168").Value = 'Agree/Accept'
$ws.Range("I169").Value = 'aa'
$ws.Range("J169").Value = 'Agree/Accept'
$ws.Range("I171").Value = 'aa'
$ws.Range("J171").Value = 'Agree/Accept'
$ws.Range("I178").Value = 'sv'
$ws.Range("J178").Value = 'Statement-opinion'
$ws.Range("I187").Value = 'aa'
$ws.Range("J187").Value = 'Agree/Accept'
$ws.Range("I192").Value = 'aa'
$ws.Range("J192").Value = 'Agree/Accept'
$ws.Range("I200").Value = 'sv'
$ws.Range("J200").Value = 'Statement-opinion'
$ws.Range("I206").Value = 'sd'
$ws.Range("J206").Value = 'Statement-non-opinion'
$ws.Range("I207").Value = 'sd'
$ws.Range("J207").Value = 'Statement-non-opinion'
$ws.Range("I216").Value = 'aa'
$ws.Range("J216").Value = 'Agree/Accept'
$ws.Range("I226").Value = 'aa'
$ws.Range("J226").Value = 'Agree/Accept'
$ws.Range("I229").Value = 'aa'
$ws.Range("J229").Value = 'Agree/Accept'
$ws.Range("I240").Value = 'aa'
$ws.Range("J240").Value = 'Agree/Accept'
$ws.Range("I244").Value = 'sv'
$ws.Range("J244").Value = 'Statement-opinion'
$ws.Range("I251").Value = 'aa'
$ws.Range("J251").Value = 'Agree/Accept'
$ws.Range("I259").Value = 'sd'
$ws.Range("J259").Value = 'Statement-non-opinion'
$ws.Range("I260").Value = 'sd'
$ws.Range("J260").Value = 'Statement-non-opinion'
$ws.Range("I263").Value = 'b'
$ws.Range("J263").Value = 'Acknowledge (Backchannel)'
$ws.Range("I266").Value = 'ba'
$ws.Range("J266").Value = 'Appreciation'
$ws.Range("I267").Value = 'sd'
$ws.Range("J267").Value = 'Statement-non-opinion'
$ws.Range("I269").Value = 'sd'
$ws.Range("J269").Value = 'Statement-non-opinion'
$ws.Range("I270").Value = '%'
$ws.Range("J270").Value = 'Uninterpretable'
$ws.Range("I280").Value = 'sv'
$ws.Range("J280").Value = 'Statement-opinion'
$ws.Range("I284").Value = 'aa'
$ws.Range("J284").Value = 'Agree/Accept'
$ws.Range("I285").Value = 'aa'
$ws.Range("J285").Value = 'Agree/Accept'
$ws.Range("I287").Value = 'sv'
$ws.Range("J287").Value = 'Statement-opinion'
$ws.Range("I297").Value = 'sv'
$ws.Range("J297").Value = 'Statement-opinion'
$ws.Range("I299").Value = 'qy'
$ws.Range("J299").Value = 'Yes-No-Question'
$ws.Range("I312").Value = 'sd'
$ws.Range("J312").Value = 'Statement-non-opinion'
$ws.Range("I315").Value = '%'
$ws.Range("J315").Value = 'Uninterpretable'
$ws.Range("I316").Value = 'b'
$ws.Range("J316").Value = 'Acknowledge (Backchannel)'
$ws.Range("I325").Value = 'aa'
$ws.Range("J325").Value = 'Agree/Accept'
$ws.Range("I337").Value = 'sv'
$ws.Range("J337").Value = 'Statement-opinion'
$ws.Range("I347").Value = 'sd'
$ws.Range("J347").Value = 'Statement-non-opinion'
$ws.Range("I350").Value = '%'
$ws.Range("J350").Value = 'Uninterpretable'
$ws.Range("I356").Value = 'sd'
$ws.Range("J356").Value = 'Statement-non-opinion'
$ws.Range("I363").Value = 'aa'
$ws.Range("J363").Value = 'Agree/Accept'
$ws.Range("I364").Value = 'aa'
$ws.Range("J364").Value = 'Agree/Accept'
$ws.Range("I365").Value = 'sd'
$ws.Range("J365").Value = 'Statement-non-opinion'
$ws.Range("I382").Value = 'b'
$ws.Range("J382").Value = 'Acknowledge (Backchannel)'
$ws.Range("I394").Value = 'aa'
$ws.Range("J394").Value = 'Agree/Accept'
$ws.Range("I400").Value = 'b'
$ws.Range("J400").Value = 'Acknowledge (Backchannel)'
$ws.Range("I403").Value = 'aa'
$ws.Range("J403").Value = 'Agree/Accept'
$ws.Range("I416").Value = 'sd'
$ws.Range("J416").Value = 'Statement-non-opinion'
$ws.Range("I421").Value = 'b'
$ws.Range("J421").Value = 'Acknowledge (Backchannel)'
$ws.Range("I423").Value = 'sd'
$ws.Range("J423").Value = 'Statement-non-opinion'
$ws.Range("I430").Value = 'aa'
$ws.Range("J430").Value = 'Agree/Accept'
$ws.Range("I431").Value = 'aa'
$ws.Range("J431").Value = 'Agree/Accept'
$ws.Range("I434").Value = 'aa'
$ws.Range("J434").Value = 'Agree/Accept'
$ws.Range("I435").Value = 'aa'
$ws.Range("J435").Value = 'Agree/Accept'
$ws.Range("I437").Value = 'ba'
$ws.Range("J437").Value = 'Appreciation'
$ws.Range("I439").Value = 'b'
$ws.Range("J439").Value = 'Acknowledge (Backchannel)'
